$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Record the previous (V3) frame numbers for rows 52-55 into the new
#    column E (before column B is overwritten with the improved frame data),
#    and add the new "diff vs previous" formulas in column F.
# ---------------------------------------------------------------------------

# Row 52 - new E/F pair (B52 did not change, this is a brand new comparison row)
$ws.Range("B52").Copy($ws.Range("E52")) | Out-Null
$ws.Range("E52").Value = 120608
$ws.Range("D52").Copy($ws.Range("F52")) | Out-Null
$ws.Range("F52").Formula = '=IF(E52<>"",IF(B52<>"",E52-B52,"-"), "-")'

# Row 53 - E53 keeps the old B53 value (121181), B53 gets the improved frame
$ws.Range("B53").Copy($ws.Range("E53")) | Out-Null
$ws.Range("E53").Value = 121181
$ws.Range("D52").Copy($ws.Range("F53")) | Out-Null
$ws.Range("F53").Formula = '=IF(E53<>"",IF(B53<>"",E53-B53,"-"), "-")'
$ws.Range("B53").Value = 121103

# Row 54 - E54 keeps the old B54 value (122496), B54 gets the improved frame
$ws.Range("B54").Copy($ws.Range("E54")) | Out-Null
$ws.Range("E54").Value = 122496
$ws.Range("D52").Copy($ws.Range("F54")) | Out-Null
$ws.Range("F54").Formula = '=IF(E54<>"",IF(B54<>"",E54-B54,"-"), "-")'
$ws.Range("B54").Value = 122418

# Row 55 - E55 keeps the old B55 value (122755), B55 gets the improved frame
$ws.Range("B55").Copy($ws.Range("E55")) | Out-Null
$ws.Range("E55").Value = 122755
$ws.Range("D52").Copy($ws.Range("F55")) | Out-Null
$ws.Range("F55").Formula = '=IF(E55<>"",IF(B55<>"",E55-B55,"-"), "-")'
$ws.Range("B55").Value = 122676

# ---------------------------------------------------------------------------
# 2) Rows that only gained a wider used-range (columns E:F now considered
#    part of the row) but received no new data of their own.
# ---------------------------------------------------------------------------
foreach ($r in 49,50,51,56,57,58,59,60,61,62,63,64) {
    $cell = $ws.Range("F" + $r)
    $cell.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# 3) Update the view: freeze pane still only covers the header row, but the
#    visible area is scrolled down and the active selection moved.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select() | Out-Null
$win.FreezePanes = $true
$win.ScrollRow = 34
$win.ScrollColumn = 1
$ws.Range("B56").Select() | Out-Null
